# FullAutoavaluacio_IsmaelSanchez.xlsx - "Add files via upload" re-edit
#
# The student re-submitted the self-assessment sheet: two of the "did you
# work on this topic" checkboxes (L19 / L30, last column "Loop testing") got
# flipped from 1 to 0, which ripples through the hidden weighting rows
# (L20/L31 -> 0, and the practice-totals C22/C33 drop from 0.9 to 0.85), and
# the final grade the student is requesting (C37) got filled in with 5,
# which flows into C39 (C37*C35). The active cell also moved to C10 (the
# name field) by the time the file was saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19 table ("Marca amb un 1 ... has treballat tú ...") ---------------
# Last topic column (L, "Loop testing") flips from 1 to 0.
$ws.Range("L19").Value = 0

# --- Row 30 table ("Marca amb un 1 ... heu treballat ...") -----------------
# Same column, same flip.
$ws.Range("L30").Value = 0

# --- Requested grade ---------------------------------------------------------
# "Nota màxima a la que optes:" (B37) now has an answer filled in to its
# right, which multiplies through to the "Total alumne:" cell (C39).
$ws.Range("C37").Value = 5

# --- Selection left on the name field (C10) when the file was saved --------
$ws.Range("C10").Select()
